# Res_sys_sim_diagram.pptx - "Mass balance function corrected"
#
#   Qreg_rel[t] = np.array([Qreg_rel[t], s[t] - s_min + I[t] + Qreg_inf[t] - E[t] - env[t]]).max()
#
# was replaced by
#
#   Qreg_rel[t] = np.array([Qreg_rel[t], s[t] - s_min + I[t] + Qreg_inf[t] - E[t] - env[t]]).min()
#
# which, on the slide itself, shows up as the formula textbox changing
# "= max(Q" to "= min(Q". The deck was also re-saved on a later day, so
# every "last updated" date field on the slide master / layouts rolled
# from 14/01/2020 to 20/01/2020.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Formula textbox: "...= max(Qreg_rel , ..." -> "...= min(Qreg_rel , ..."
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$formulaShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text.IndexOf("= max(Q") -ge 0) {
            $formulaShape = $sh
            break
        }
    }
}

if ($formulaShape -ne $null) {
    $tr = $formulaShape.TextFrame.TextRange
    $old = "= max(Q"
    $new = "= min(Q"
    $pos = $tr.Text.IndexOf($old)
    if ($pos -ge 0) {
        $run = $tr.Characters($pos + 1, $old.Length)
        $run.Text = $new
    }
}

# ---------------------------------------------------------------------
# 2) "Update automatically" date placeholder on the master + every
#    slide layout: 14/01/2020 -> 20/01/2020
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes, [string]$oldDate, [string]$newDate) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "14/01/2020" "20/01/2020"

for ($k = 1; $k -le $master.CustomLayouts.Count; $k++) {
    $layout = $master.CustomLayouts.Item($k)
    Update-DatePlaceholders $layout.Shapes "14/01/2020" "20/01/2020"
}
